# idRuntimeException-expected-generation.docx was regenerated against a
# newer build (3.1.0 -> 3.1.1) of the Acceleo/M2Doc/JUnit stack, which
# shifted several source line numbers and stack-trace frames in the
# Java exception text embedded in the second paragraph of the document.
# Replace the affected block of the stack trace (bounded by unchanged
# context lines on both sides, so the match is unambiguous even though
# several frames repeat verbatim elsewhere in the trace).

$d = $word.ActiveDocument

$searchText = "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)`n`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1275)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:160)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:281)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:805)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)`n`tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)`n`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)`n`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)"
$replaceText = "`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:186)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)`n`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1278)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:160)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1467)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:297)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:282)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:845)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)`n`tat sun.reflect.GeneratedMethodAccessor73.invoke(Unknown Source)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)`n`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)`n`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)"

$range = $d.Content
$found = $range.Find.Execute($searchText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)

if (-not $found) {
    throw "idRuntimeException: stack trace block to replace was not found."
}
